# New trade row observed 20 minutes after the previous one (row 3) - append it
# to the trade log, copying row 3's cell formatting (date format on A/G) first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:I3").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 42641.546990740739
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 9962.4699999999993
$ws.Range("D4").Value = 9955.5
$ws.Range("E4").Value = 106.51
$ws.Range("F4").Value = 106.36
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = -0.14000000000000001
$ws.Range("I4").Value = $false
